# Updates cryptos list prices / 1h-volume changes (and two coin-name /
# link swaps in rows 39-40 and the row-51 coin replacement), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Most Price (column D) values look like plain decimals (e.g. "1.000",
# "238.82"), so a direct .Value assignment would get silently coerced by
# Excel into a number (dropping trailing zeros / losing exact text).
# To keep them as literal text - as they were authored - we briefly force
# the cell to Text format, assign the string, then reset the cell style
# back to "Normal" so no stray style index is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.396.37"

$ws.Range("D3").Value = "1.875.80"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4803"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2817"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("D10").Value = "1.872.13"
$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07491"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.070"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.45"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6626"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").Value = "30.351.27"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007588"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.06%  "

$ws.Range("D20").Value = "2.114.41"
$ws.Range("E20").Value = "  -1.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.301"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.89%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "219.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.189"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.332"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.960"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("E29").Value = "  +0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09380"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.307"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.025"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05016"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.208"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7427"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.705"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01827"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.617"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.19%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9052"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.059"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.853"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4271"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.426"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1273"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.474"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.881"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.53%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3881"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "
